$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.471.95'
$ws.Range('E2').Value = '  +1.11%  '
$ws.Range('D3').Value = '1.828.89'
$ws.Range('E3').Value = '  +1.16%  '
$ws.Range('E4').Value = '  +0.23%  '
$ws.Range('D5').Value = "'315.13"
$ws.Range('E5').Value = '  -0.59%  '
$ws.Range('E6').Value = '  +0.13%  '
$ws.Range('D7').Value = "'0.5121"
$ws.Range('E7').Value = '  -3.46%  '
$ws.Range('D8').Value = "'0.3916"
$ws.Range('E8').Value = '  -1.86%  '
$ws.Range('D9').Value = "'0.07676"
$ws.Range('E9').Value = '  +1.55%  '
$ws.Range('D10').Value = "'41.82"
$ws.Range('E10').Value = '  +0.46%  '
$ws.Range('D11').Value = "'1.112"
$ws.Range('E11').Value = '  +1.81%  '
$ws.Range('D12').Value = "'21.08"
$ws.Range('E12').Value = '  +2.82%  '
$ws.Range('D13').Value = "'6.283"
$ws.Range('E13').Value = '  +0.45%  '
$ws.Range('E14').Value = '  +0.26%  '
$ws.Range('D15').Value = "'7.537"
$ws.Range('E15').Value = '  -0.52%  '
$ws.Range('D16').Value = '1.823.19'
$ws.Range('E16').Value = '  +1.24%  '
$ws.Range('D17').Value = "'93.73"
$ws.Range('E17').Value = '  +5.29%  '
$ws.Range('D18').Value = "'0.00001102"
$ws.Range('E18').Value = '  +3.41%  '
$ws.Range('D19').Value = "'0.06722"
$ws.Range('E19').Value = '  +2.15%  '
$ws.Range('E20').Value = '  +2.23%  '
$ws.Range('E21').Value = '  +0.14%  '
$ws.Range('D22').Value = "'6.163"
$ws.Range('E22').Value = '  +2.33%  '
$ws.Range('D23').Value = '28.500.91'
$ws.Range('E23').Value = '  +1.21%  '
$ws.Range('D24').Value = "'11.17"
$ws.Range('E24').Value = '  +0.34%  '
$ws.Range('D25').Value = "'2.256"
$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').Value = "'20.68"
$ws.Range('E26').Value = '  +1.77%  '
$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').Value = "'156.61"
$ws.Range('E27').Value = '  +0.03%  '
$ws.Range('D28').Value = '2.038.09'
$ws.Range('E28').Value = '  +1.31%  '
$ws.Range('D29').Value = "'2.395"
$ws.Range('E29').Value = '  +0.49%  '
$ws.Range('D30').Value = "'124.42"
$ws.Range('E30').Value = '  +1.20%  '
$ws.Range('D31').Value = "'1.116"
$ws.Range('E31').Value = '  +1.66%  '
$ws.Range('D32').Value = "'0.1090"
$ws.Range('E32').Value = '  -0.61%  '
$ws.Range('D33').Value = "'5.668"
$ws.Range('E33').Value = '  +2.14%  '
$ws.Range('D34').Value = "'3.655"
$ws.Range('E34').Value = '  -0.43%  '
$ws.Range('D35').Value = "'0.07030"
$ws.Range('E35').Value = '  -0.82%  '
$ws.Range('D36').Value = "'0.2215"
$ws.Range('E36').Value = '  -0.27%  '
$ws.Range('D37').Value = "'8.968"
$ws.Range('E37').Value = '  +4.54%  '
$ws.Range('D38').Value = "'0.02324"
$ws.Range('E38').Value = '  +0.98%  '
$ws.Range('D39').Value = "'5.159"
$ws.Range('E39').Value = '  -1.13%  '
$ws.Range('D40').Value = "'0.6266"
$ws.Range('E40').Value = '  +1.25%  '
$ws.Range('D41').Value = "'11.22"
$ws.Range('E41').Value = '  -0.45%  '
$ws.Range('D42').Value = "'1.182"
$ws.Range('E42').Value = '  -1.24%  '
$ws.Range('E43').Value = '  +0.08%  '
$ws.Range('D44').Value = "'1.392"
$ws.Range('E44').Value = '  -0.84%  '
$ws.Range('D45').Value = "'13.43"
$ws.Range('E45').Value = '  +0.32%  '
$ws.Range('D46').Value = "'0.5904"
$ws.Range('E46').Value = '  +2.47%  '
$ws.Range('D47').Value = "'3.716"
$ws.Range('E47').Value = '  +0.86%  '
$ws.Range('D48').Value = "'125.00"
$ws.Range('E48').Value = '  -0.09%  '
$ws.Range('D49').Value = "'1.981"
$ws.Range('E49').Value = '  +2.19%  '
$ws.Range('D50').Value = "'1.198"
$ws.Range('E50').Value = '  +1.06%  '
$ws.Range('E51').Value = '  +1.55%  '
